$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "64.183.29"
$ws.Range("E2").Value2 = "  -0.79%  "
$ws.Range("D3").Value2 = "3.399.97"
$ws.Range("E3").Value2 = "  -1.19%  "
$ws.Range("E4").Value2 = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "569.72"
$ws.Range("E5").Value2 = "  -0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "155.62"
$ws.Range("E6").Value2 = "  -2.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.623"
$ws.Range("E7").Value2 = "  +6.97%  "
$ws.Range("E8").Value2 = "  +0.06%  "
$ws.Range("D9").Value2 = "3.400.68"
$ws.Range("E9").Value2 = "  -1.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "7.13"
$ws.Range("E10").Value2 = "  -3.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.122"
$ws.Range("E11").Value2 = "  -2.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.440"
$ws.Range("E12").Value2 = "  -0.14%  "
$ws.Range("D13").Value2 = "3.984.91"
$ws.Range("E13").Value2 = "  -1.31%  "
$ws.Range("E14").Value2 = "  +0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.0000186"
$ws.Range("E15").Value2 = "  -3.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "27.64"
$ws.Range("E16").Value2 = "  -2.04%  "
$ws.Range("D17").Value2 = "64.194.27"
$ws.Range("E17").Value2 = "  -0.84%  "
$ws.Range("D18").Value2 = "3.434.32"
$ws.Range("E18").Value2 = "  -1.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "6.30"
$ws.Range("E19").Value2 = "  -0.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "13.86"
$ws.Range("E20").Value2 = "  -2.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "376.23"
$ws.Range("E21").Value2 = "  -2.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "8.00"
$ws.Range("E22").Value2 = "  -2.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "0.545"
$ws.Range("E23").Value2 = "  +0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "0.999"
$ws.Range("E24").Value2 = "  -0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "71.93"
$ws.Range("E25").Value2 = "  -1.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "0.0000118"
$ws.Range("E26").Value2 = "  -4.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "10.33"
$ws.Range("E27").Value2 = "  +6.20%  "
$ws.Range("E28").Value2 = "  -1.68%  "
$ws.Range("E29").Value2 = "  +0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "1.48"
$ws.Range("E30").Value2 = "  +2.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "6.15"
$ws.Range("E31").Value2 = "  -0.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "2.00"
$ws.Range("E32").Value2 = "  -2.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "23.00"
$ws.Range("E33").Value2 = "  -2.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "7.14"
$ws.Range("E34").Value2 = "  +0.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "1.60"
$ws.Range("E35").Value2 = "  +6.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "159.48"
$ws.Range("E36").Value2 = "  -2.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "1.89"
$ws.Range("E37").Value2 = "  -0.63%  "
$ws.Range("B38").Value2 = "Hedera"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.0761"
$ws.Range("E38").Value2 = "  -0.51%  "
$ws.Range("B39").Value2 = "RenderToken"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "6.86"
$ws.Range("E39").Value2 = "  +5.03%  "
$ws.Range("D40").Value2 = "2.875.55"
$ws.Range("E40").Value2 = "  -4.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "26.47"
$ws.Range("E41").Value2 = "  -2.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "4.60"
$ws.Range("E42").Value2 = "  +0.18%  "
$ws.Range("B43").Value2 = "OKB"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "42.76"
$ws.Range("E43").Value2 = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.0315"
$ws.Range("E44").Value2 = "  +0.00%  "
$ws.Range("B45").Value2 = "InjectiveProtocol"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "26.25"
$ws.Range("E45").Value2 = "  +6.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.767"
$ws.Range("E46").Value2 = "  -0.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "319.64"
$ws.Range("E47").Value2 = "  +5.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "1.08"
$ws.Range("E48").Value2 = "  -1.22%  "
$ws.Range("E49").Value2 = "  +2.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "2.17"
$ws.Range("E50").Value2 = "  -0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.858"
$ws.Range("E51").Value2 = "  -2.14%  "
